$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "52.992.01"
$ws.Range("E2").Value = "  -10.28%  "
$ws.Range("D3").Value = "2.375.56"
$ws.Range("E3").Value = "  -13.06%  "
$ws.Range("E4").Value = "  +0.19%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "455.00"
$c.Style = "Normal"
$ws.Range("E5").Value = "  -10.55%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "127.65"
$c.Style = "Normal"
$ws.Range("E6").Value = "  -10.15%  "
$c = $ws.Range("D7")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E7").Value = "  +0.02%  "
$c = $ws.Range("D8")
$c.NumberFormat = "@"
$c.Value = "0.478"
$c.Style = "Normal"
$ws.Range("E8").Value = "  -10.67%  "
$ws.Range("D9").Value = "2.390.85"
$ws.Range("E9").Value = "  -13.00%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.0936"
$c.Style = "Normal"
$ws.Range("E10").Value = "  -10.75%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "5.25"
$c.Style = "Normal"
$ws.Range("E11").Value = "  -14.23%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.308"
$c.Style = "Normal"
$ws.Range("E12").Value = "  -11.80%  "
$ws.Range("E13").Value = "  -4.71%  "
$ws.Range("D14").Value = "2.791.77"
$ws.Range("E14").Value = "  -12.94%  "
$ws.Range("D15").Value = "53.112.31"
$ws.Range("E15").Value = "  -9.93%  "
$c = $ws.Range("D16")
$c.NumberFormat = "@"
$c.Value = "19.34"
$c.Style = "Normal"
$ws.Range("E16").Value = "  -11.61%  "
$ws.Range("E17").Value = "  -5.73%  "
$ws.Range("D18").Value = "2.392.85"
$ws.Range("E18").Value = "  -12.76%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "4.10"
$c.Style = "Normal"
$ws.Range("E19").Value = "  -13.84%  "
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "302.80"
$c.Style = "Normal"
$ws.Range("E20").Value = "  -12.49%  "
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "9.22"
$c.Style = "Normal"
$ws.Range("E21").Value = "  -16.36%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "0.996"
$c.Style = "Normal"
$ws.Range("E22").Value = "  -0.23%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "5.67"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +0.99%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "5.26"
$c.Style = "Normal"
$ws.Range("E24").Value = "  -16.13%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "55.31"
$c.Style = "Normal"
$ws.Range("E25").Value = "  -12.56%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +0.30%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "0.379"
$c.Style = "Normal"
$ws.Range("E27").Value = "  -11.30%  "
$ws.Range("B28").Value = "WrappedeETH"
$ws.Range("C28").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D28").Value = "2.468.23"
$ws.Range("E28").Value = "  -13.32%  "
$ws.Range("B29").Value = "Kaspa"
$ws.Range("C29").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "0.150"
$c.Style = "Normal"
$ws.Range("E29").Value = "  -13.56%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "6.98"
$c.Style = "Normal"
$ws.Range("E30").Value = "  -7.25%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "0.998"
$c.Style = "Normal"
$ws.Range("E31").Value = "  -0.04%  "
$ws.Range("D32").Value = "0.0₃0708"
$ws.Range("E32").Value = "  -16.11%  "
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "144.72"
$c.Style = "Normal"
$ws.Range("E33").Value = "  -3.21%  "
$c = $ws.Range("D34")
$c.NumberFormat = "@"
$c.Value = "17.44"
$c.Style = "Normal"
$ws.Range("E34").Value = "  -9.20%  "
$c = $ws.Range("D35")
$c.NumberFormat = "@"
$c.Value = "1.40"
$c.Style = "Normal"
$ws.Range("E35").Value = "  -13.70%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "4.92"
$c.Style = "Normal"
$ws.Range("E36").Value = "  -9.35%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "3.46"
$c.Style = "Normal"
$ws.Range("E37").Value = "  -17.94%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "1.04"
$c.Style = "Normal"
$ws.Range("E38").Value = "  -8.88%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "0.783"
$c.Style = "Normal"
$ws.Range("E39").Value = "  -18.31%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "0.995"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +0.10%  "
$c = $ws.Range("D41")
$c.NumberFormat = "@"
$c.Value = "32.70"
$c.Style = "Normal"
$ws.Range("E41").Value = "  -9.55%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "0.586"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -3.55%  "
$ws.Range("B43").Value = "Hedera"
$ws.Range("C43").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$c = $ws.Range("D43")
$c.NumberFormat = "@"
$c.Value = "0.0518"
$c.Style = "Normal"
$ws.Range("E43").Value = "  -7.54%  "
$ws.Range("B44").Value = "Filecoin"
$ws.Range("C44").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "3.22"
$c.Style = "Normal"
$ws.Range("E44").Value = "  -9.11%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "10.12"
$c.Style = "Normal"
$ws.Range("E45").Value = "  -2.30%  "
$ws.Range("E46").Value = "  -12.97%  "
$ws.Range("D47").Value = "1.928.05"
$ws.Range("E47").Value = "  -11.79%  "
$c = $ws.Range("D48")
$c.NumberFormat = "@"
$c.Value = "0.0214"
$c.Style = "Normal"
$ws.Range("E48").Value = "  -6.19%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "0.0855"
$c.Style = "Normal"
$ws.Range("E49").Value = "  -3.72%  "
$c = $ws.Range("D50")
$c.NumberFormat = "@"
$c.Value = "4.13"
$c.Style = "Normal"
$ws.Range("E50").Value = "  -14.17%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "16.22"
$c.Style = "Normal"
$ws.Range("E51").Value = "  -15.17%  "
